$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Lab")

# --- Fix existing translation texts (PG/VG -> VG/PG wording + shortened quick-mix subtitle) ---
$found = $ws.Range("B1:B56").Find("lab.liquid.list.empty.subtitle", [System.Type]::Missing, [System.Type]::Missing, 1)
$found.Offset(0, 1).Value = "Přidám liquidu si jej zaevidujete a můžete sledovat různé zajímavé hodnoty, jako jsou například datum jeho dozrání, datum expirace (pokud to stihne) nebo vypočtené poměry VG/PG."

$found = $ws.Range("B1:B56").Find("lab.liquid.preview.aroma.pgvg", [System.Type]::Missing, [System.Type]::Missing, 1)
$found.Offset(0, 1).Value = "Poměr VG/PG aromatu"

$found = $ws.Range("B1:B56").Find("lab.liquid.preview.base.pgvg", [System.Type]::Missing, [System.Type]::Missing, 1)
$found.Offset(0, 1).Value = "Poměr VG/PG báze"

$found = $ws.Range("B1:B56").Find("lab.liquid.create.quick-mix.subtitle", [System.Type]::Missing, [System.Type]::Missing, 1)
$found.Offset(0, 1).Value = "Tato funkce umožní namíchání liquidu formou dolití bází, kdy už víte, co děláte a není třeba složitého vyklikávání dalších údajů."

# --- Append new translation rows for lab.liquid.hint.vg.* keys ---
$lastRow = $ws.Cells.Item(1, 1).SpecialCells(11).Row
$lastRow = $lastRow + 1
$ws.Range("A2:C2").Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($lastRow, 1).Value = "cs"
$ws.Cells.Item($lastRow, 2).Value = "lab.liquid.hint.vg.30"
$ws.Cells.Item($lastRow, 3).Value = "Liquid vhodný pro MTL, může velmi dobře nést chuť, nicméně díky velkému poměru PG bude velmi řídký, tudíž s tím bude třeba počítat při buildu nebo použití v podu."

$lastRow = $lastRow + 1
$ws.Range("A2:C2").Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($lastRow, 1).Value = "cs"
$ws.Cells.Item($lastRow, 2).Value = "lab.liquid.hint.vg.40"
$ws.Cells.Item($lastRow, 3).Value = "Tento poměr bude spíše tekutější (díky většímu poměru PG) a bude lépe nést chuť. Vhodnější do MTL s patřičným buildem (může protékat)."

$lastRow = $lastRow + 1
$ws.Range("A2:C2").Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($lastRow, 1).Value = "cs"
$ws.Cells.Item($lastRow, 2).Value = "lab.liquid.hint.vg.50"
$ws.Cells.Item($lastRow, 3).Value = "Běžný poměr, který je relativně řídký, tudíž je vhodný do MTL (díky PG nese chuť), je ovšem třeba mít pro tento poměr správný build, protože může protékat."

$lastRow = $lastRow + 1
$ws.Range("A2:C2").Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($lastRow, 1).Value = "cs"
$ws.Cells.Item($lastRow, 2).Value = "lab.liquid.hint.vg.60"
$ws.Cells.Item($lastRow, 3).Value = "Poměrně běžný poměr který je spíše vhodný pro MTL vaping (tudíž je třeba mít patřičný build), případně pody, které tento poměr snesou. Celkový mix bude mírně hustější (díky trochu většímu množství VG)."

$lastRow = $lastRow + 1
$ws.Range("A2:C2").Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($lastRow, 1).Value = "cs"
$ws.Cells.Item($lastRow, 2).Value = "lab.liquid.hint.vg.70"
$ws.Cells.Item($lastRow, 3).Value = "Celkem oblíbený poměr spíše určený pro DL vaping, jelikož výsledná směs bude produkovat více páry na úkor chuti. Také se jedná o hustější liquid, tudíž je potřeba s tímto počítat u buildu."

$lastRow = $lastRow + 1
$ws.Range("A2:C2").Copy()
$ws.Cells.Item($lastRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($lastRow, 1).Value = "cs"
$ws.Cells.Item($lastRow, 2).Value = "lab.liquid.hint.vg.80"
$ws.Cells.Item($lastRow, 3).Value = "Jedná se o poměrně extrémní poměr více hrčen pro DL vapování, kde je potlačena chuť ve prospěch velkých mraků. Celkový mix bude hustý, proto je třeba s tímto faktem počítat v buildech (hustý liquid nemusí být schopný rychle zásobit vatu u spirálky)."

# --- Sort A2:C<lastRow> ascending by column B (translation key) ---
$sortRange = $ws.Range("A2:C" + $lastRow)
$sortKey = $ws.Range("B2:B" + $lastRow)
$sortRange.Sort($sortKey, 1)

